# Generate Report for Handback
# Replaces the first handback file's GUID-named artifacts
# (87125e70-4812-4d7a-bab9-591f8a17caf5 -> 33d99b1b-f6df-4c6b-946b-effcbb21a229)
# and the second handback file's GUID-named artifacts
# (f9fde33b-ca70-47b4-998d-c05cc45437ce -> ffff809b08e7-ae90-45ee-93b3-582579312b3e)
# with refreshed handoff/handback timestamps and xliff content hashes.

$wb = $excel.ActiveWorkbook

$oldGuid1 = "87125e70-4812-4d7a-bab9-591f8a17caf5"
$newGuid1 = "33d99b1b-f6df-4c6b-946b-effcbb21a229"
$oldGuid2 = "f9fde33b-ca70-47b4-998d-c05cc45437ce"
$newGuid2 = "ffff809b08e7-ae90-45ee-93b3-582579312b3e"

$oldHash1 = "cd1ba9ed196f88b77849019ca7d56c83231924d8"
$newHash1 = "660067b065ccdcb6ad2a13e5b5e2f7606ee8ed26"

$newZhCnXlf = "$newGuid1.$newHash1.zh-cn.xlf"
$newDeDeXlf = "$newGuid1.$newHash1.de-de.xlf"

$newOverviewDate = "2016-08-13 17:28:55"
$newZhCnHandoffDate = "2016-08-13 17:28:47"
$newZhCnHandbackDate = "2016-08-13 17:29:16"
$newDeDeDate = "2016-08-13 17:29:26"

$repoAddr = "https://github.com/OpenLocalizationTestOrg/oltest/blob/70c715776eb314f56c976c5624cd64c598547645/e2e"
$zhcnAddr = "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/3cc1afd553554e995bb645b59de90e79d7be7775/e2e"
$dedeAddr = "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/80a7e5d20a7b918046216419b6b76605caa8c2b3/e2e"

# Cornflower-blue link color (matches the workbook's original custom "HyperLink"
# cell style, RGB 6495ED) expressed as the BGR integer Excel's Font.Color expects.
$linkColor = 15570276

function Restyle-Hyperlink($range) {
    $range.Font.Color = $linkColor
    $range.Font.Underline = 2
}

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = "$newGuid1.md"
$wsOverview.Range("B2").Value = "e2e\$newGuid1.md"
$wsOverview.Range("G2").Value = $newOverviewDate

$wsOverview.Range("A3").Value = "$newGuid2.md"
$wsOverview.Range("B3").Value = "e2e\$newGuid2.md"
$wsOverview.Range("G3").Value = $newOverviewDate

$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), "$repoAddr/$newGuid1.md", "", "", "e2e\$newGuid1.md")
$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), "$repoAddr/$newGuid2.md", "", "", "e2e\$newGuid2.md")
Restyle-Hyperlink $wsOverview.Range("B2")
Restyle-Hyperlink $wsOverview.Range("B3")

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("A2").Value = "$newGuid1.md"
$wsZhCn.Range("G2").Value = $newZhCnXlf
$wsZhCn.Range("H2").Value = $newZhCnHandoffDate
$wsZhCn.Range("I2").Value = "$newGuid1.md"
$wsZhCn.Range("J2").Value = $newZhCnXlf
$wsZhCn.Range("K2").Value = $newZhCnHandbackDate

$wsZhCn.Range("A3").Value = "$newGuid2.md"
$wsZhCn.Range("G3").Value = $newZhCnXlf
$wsZhCn.Range("H3").Value = $newZhCnHandoffDate
$wsZhCn.Range("I3").Value = "$newGuid2.md"
$wsZhCn.Range("J3").Value = $newZhCnXlf
$wsZhCn.Range("K3").Value = $newZhCnHandbackDate

$wsZhCn.Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), "$repoAddr/$newGuid1.md", "", "", "$newGuid1.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I2"), "$zhcnAddr/$newGuid1.md", "", "", "$newGuid1.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), "$repoAddr/$newGuid2.md", "", "", "$newGuid2.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I3"), "$zhcnAddr/$newGuid2.md", "", "", "$newGuid2.md")
Restyle-Hyperlink $wsZhCn.Range("A2")
Restyle-Hyperlink $wsZhCn.Range("I2")
Restyle-Hyperlink $wsZhCn.Range("A3")
Restyle-Hyperlink $wsZhCn.Range("I3")

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("A2").Value = "$newGuid1.md"
$wsDeDe.Range("G2").Value = $newDeDeXlf
$wsDeDe.Range("H2").Value = $newOverviewDate
$wsDeDe.Range("I2").Value = "$newGuid1.md"
$wsDeDe.Range("J2").Value = $newDeDeXlf
$wsDeDe.Range("K2").Value = $newDeDeDate

$wsDeDe.Range("A3").Value = "$newGuid2.md"
$wsDeDe.Range("G3").Value = $newDeDeXlf
$wsDeDe.Range("H3").Value = $newOverviewDate
$wsDeDe.Range("I3").Value = "$newGuid2.md"
$wsDeDe.Range("J3").Value = $newDeDeXlf
$wsDeDe.Range("K3").Value = $newDeDeDate

$wsDeDe.Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), "$repoAddr/$newGuid1.md", "", "", "$newGuid1.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I2"), "$dedeAddr/$newGuid1.md", "", "", "$newGuid1.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), "$repoAddr/$newGuid2.md", "", "", "$newGuid2.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I3"), "$dedeAddr/$newGuid2.md", "", "", "$newGuid2.md")
Restyle-Hyperlink $wsDeDe.Range("A2")
Restyle-Hyperlink $wsDeDe.Range("I2")
Restyle-Hyperlink $wsDeDe.Range("A3")
Restyle-Hyperlink $wsDeDe.Range("I3")
